$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns are treated as text so
# formatted numeric-looking strings (e.g. "1.000", "0.9950") are preserved exactly.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.160.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.73%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9950"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.65"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.25%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3922"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3899"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.81%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.72"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.379"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9948"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08517"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.69%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.25"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.297"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.156"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +8.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001324"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.28%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.648.41"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.09"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06968"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.82%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.21"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.976"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.77"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.122.06"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.174"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +10.99%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.505"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.33"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.94"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.75"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.339"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.28%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.919"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.487"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.832.79"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.27%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.060"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +10.26%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03048"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.98%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08198"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.92%  "

# Row 37
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.37"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +10.76%  "

# Row 38
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.780"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.27%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2742"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.10%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09213"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.37%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7636"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.63%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.61"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.50%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.428"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.40%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.76%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7020"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.89%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.517"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.105"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.08%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08349"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.39"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.35%  "

# Row 51
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.253"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.86%  "
